$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new rows 22-30 (regcntr_id, usr_id increment from existing data)
$data = @(
    @(10002, 110021),
    @(10003, 110022),
    @(10004, 110023),
    @(10005, 110024),
    @(10006, 110025),
    @(10007, 110026),
    @(10008, 110027),
    @(10009, 110028),
    @(10010, 110029)
)

$row = 22
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
    $row++
}

# Set the selection to mimic the "next empty row selected" state seen after data entry
$ws.Range("A31:XFD1048576").Select()

# Configure page setup (matches xlPortrait = 1)
$ws.PageSetup.Orientation = 1
